# "Further works on the paper"
#
# 1. On fa_keys: selection moves from A2:A42 down to a single-cell
#    selection at E41 (and it's no longer the active/visible tab).
# 2. On keys_v3: the DS_* indicator columns (C:J) get reshuffled into a
#    new order, the sheet becomes the active tab, and the selection is
#    the whole of column I (as if the user clicked the column header).

$wb = $excel.ActiveWorkbook

# ---- fa_keys: just a new selection, no longer the active sheet -------
$faKeys = $wb.Worksheets.Item("fa_keys")
[void]$faKeys.Activate()
[void]$faKeys.Range("E41").Select()

# ---- keys_v3: reorder columns C..J ------------------------------------
$keysV3 = $wb.Worksheets.Item("keys_v3")
[void]$keysV3.Activate()

# Snapshot every source column *before* writing anything, since several
# of the destinations are also sources for other columns.
$colB = $keysV3.Range("B1:B42").Value2
$colC = $keysV3.Range("C1:C42").Value2
$colD = $keysV3.Range("D1:D42").Value2
$colE = $keysV3.Range("E1:E42").Value2
$colF = $keysV3.Range("F1:F42").Value2
$colG = $keysV3.Range("G1:G42").Value2
$colH = $keysV3.Range("H1:H42").Value2
$colI = $keysV3.Range("I1:I42").Value2
$colJ = $keysV3.Range("J1:J42").Value2

# New column order (B stays put):
#   C <- old G   D <- old J   E <- old D   F <- old E
#   G <- old I   H <- old H   I <- old C   J <- old F
$keysV3.Range("B1:B42").Value2 = $colB
$keysV3.Range("C1:C42").Value2 = $colG
$keysV3.Range("D1:D42").Value2 = $colJ
$keysV3.Range("E1:E42").Value2 = $colD
$keysV3.Range("F1:F42").Value2 = $colE
$keysV3.Range("G1:G42").Value2 = $colI
$keysV3.Range("H1:H42").Value2 = $colH
$keysV3.Range("I1:I42").Value2 = $colC
$keysV3.Range("J1:J42").Value2 = $colF

# Final selection: whole of column I, like clicking its header.
[void]$keysV3.Columns.Item(9).Select()
